$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '73.257.82'
$ws.Range("E2").Value = '  -0.19%  '

$ws.Range("D3").Value = '3.983.48'
$ws.Range("E3").Value = '  -1.76%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.91%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +12.60%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.690'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.97%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.799'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.50%  '

$ws.Range("E10").Value = '  +8.10%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.30'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.93%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000338'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.03%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.74'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.46%  '

$ws.Range("D14").Value = '4.619.19'
$ws.Range("E14").Value = '  -1.74%  '

$ws.Range("D15").Value = '3.983.54'
$ws.Range("E15").Value = '  -1.69%  '

$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.26'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.70%  '

$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.34'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.30%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '21.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.41%  '

$ws.Range("D19").Value = '73.151.58'
$ws.Range("E19").Value = '  -0.14%  '

$ws.Range("E20").Value = '  -1.08%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '464.65'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.21%  '

$ws.Range("E22").Value = '  +5.64%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '96.79'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.93%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.42'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.59%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.32'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.80%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.53%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.67'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.89'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.46'
$ws.Range("D30").Style = "Normal"

$ws.Range("E31").Value = '  +2.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '14.12'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.70%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '50.08'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.18%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.131'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.90%  '

$ws.Range("E35").Value = '  +13.60%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '70.67'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.78%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '641.09'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.15%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.435'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.08%  '

$ws.Range("E39").Value = '  -0.68%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.42'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.89%  '

$ws.Range("E41").Value = '  +0.06%  '

$ws.Range("E42").Value = '  +0.11%  '

$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.27'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +39.54%  '

$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0488'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.87%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.60'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.77%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.150'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.52%  '

$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.99'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -10.35%  '

$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000300'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.47%  '

$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.44'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.11%  '

$ws.Range("B50").Value = 'Fetch.AI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.65'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.74%  '

$ws.Range("D51").Value = '2.819.02'
$ws.Range("E51").Value = '  +1.44%  '
